$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mooneye")
$ws.Activate()

# OAM DMA memory-read fix: these Mooneye test cases now pass.
$rows = @(50, 52, 64, 65, 66, 67, 68, 69, 73, 74, 76)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "Passed"
}

# Update the active window's scroll position and selection to match where
# the user ended up after re-running the tests.
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("E77").Select()
